# Updated symbol list on Mon Jan 30 21:45:02 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# crypto rows that moved. Values are entered with a leading apostrophe
# so Excel keeps them as literal text (matching the original inlineStr
# cells) instead of auto-converting to numbers/percentages, and the
# style is reset to "Normal" right after so no stray "quote prefix"
# number format sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.66%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.08%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.085"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.90%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07719"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-5.88%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.350"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.07%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.887"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-7.94%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.200"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.78%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-6.66%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9178"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.39%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-15.15%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1880"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.42%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08704"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-4.81%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03401"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-2.94%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09703"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.85%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001375"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.14%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005911"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-4.58%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.586"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.70%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3408"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-2.13%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1289"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.45%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.016"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.54%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'6.01%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'5,171.86%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04319"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.65%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-1.20%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004536"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-5.75%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'3.84%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02204"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.19%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04919"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-5.59%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007561"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.48%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009943"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.85%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1333"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.97%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.001996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.40%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008489"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-9.05%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006547"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.86%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.002999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "'-22.88%"
$ws.Range("E49").Style = "Normal"
